# "updates for the 6th may" -- add the 6-May-2018 round to the may18 sheet,
# then leave may18 as the active sheet/tab with the next empty score cell
# (B44) selected, ready for further entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("may18")

# --- New round: 6 May 2018 (row 43 header, rows 44-61 hole-by-hole data) ---
$ws.Range("A43").Value = 43226
$ws.Range("A43").NumberFormat = "d-mmm-yy"
$ws.Range("B43").Value = "Score"
$ws.Range("C43").Value = "Fairway"
$ws.Range("D43").Value = "GIR"
$ws.Range("E43").Value = "Putts"
$ws.Range("F43").Value = "Comment"

$ws.Range("A44").Value = "Hole 1"
$ws.Range("B44").Value = 4
$ws.Range("C44").Value = "S"
$ws.Range("E44").Value = 1

$ws.Range("A45").Value = "Hole 2"
$ws.Range("B45").Value = 5
$ws.Range("E45").Value = 2

$ws.Range("A46").Value = "Hole 3"
$ws.Range("B46").Value = 5
$ws.Range("C46").Value = "S"
$ws.Range("E46").Value = 2

$ws.Range("A47").Value = "Hole 4"
$ws.Range("B47").Value = 5
$ws.Range("C47").Value = "S"
$ws.Range("E47").Value = 2

$ws.Range("A48").Value = "Hole 5"
$ws.Range("B48").Value = 4
$ws.Range("E48").Value = 1

$ws.Range("A49").Value = "Hole 6"
$ws.Range("B49").Value = 6
$ws.Range("C49").Value = "L"
$ws.Range("E49").Value = 2

$ws.Range("A50").Value = "Hole 7"
$ws.Range("B50").Value = 5
$ws.Range("C50").Value = "S"
$ws.Range("E50").Value = 2

$ws.Range("A51").Value = "Hole 8"
$ws.Range("B51").Value = 3
$ws.Range("E51").Value = 2

$ws.Range("A52").Value = "Hole 9"
$ws.Range("B52").Value = 5
$ws.Range("C52").Value = "R"
$ws.Range("E52").Value = 2

$ws.Range("A53").Value = "Hole 10"
$ws.Range("B53").Value = 5
$ws.Range("C53").Value = "R"
$ws.Range("E53").Value = 2

$ws.Range("A54").Value = "Hole 11"
$ws.Range("B54").Value = 3
$ws.Range("E54").Value = 1

$ws.Range("A55").Value = "Hole 12"
$ws.Range("B55").Value = 4
$ws.Range("C55").Value = "S"
$ws.Range("E55").Value = 1

$ws.Range("A56").Value = "Hole 13"
$ws.Range("B56").Value = 5
$ws.Range("C56").Value = "L"
$ws.Range("E56").Value = 1

$ws.Range("A57").Value = "Hole 14"
$ws.Range("B57").Value = 7
$ws.Range("C57").Value = "S"
$ws.Range("E57").Value = 2

$ws.Range("A58").Value = "Hole 15"
$ws.Range("B58").Value = 3
$ws.Range("E58").Value = 1

$ws.Range("A59").Value = "Hole 16"
$ws.Range("B59").Value = 4
$ws.Range("C59").Value = "S"
$ws.Range("E59").Value = 2

$ws.Range("A60").Value = "Hole 17"
$ws.Range("B60").Value = 4
$ws.Range("C60").Value = "L"
$ws.Range("E60").Value = 2

$ws.Range("A61").Value = "Hole 18"
$ws.Range("B61").Value = 5
$ws.Range("C61").Value = "R"
$ws.Range("E61").Value = 2

# --- View state: may18 becomes the active/selected tab, with the next
# blank score cell selected (mirrors mid-entry "Score" cell for the round
# being typed in) ---
$ws.Activate()
$ws.Range("B44").Select()
